$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New block of rows (18-21) is a repeat of the original 4-row cycle (rows 2-5),
# appended at the bottom of the sheet (dimension grows from A1:I17 to A1:I21).
$data = @(
    @(1.1993322807698887, 2.3758877717639884, 2.413486364972186, -0.86585001746396684, -0.78012023365383742, 0, 8.3366194783562833, 0.86585001746396684, 4.5494284741316866),
    @([double]"2.688821387764051e-17", 2.7925609058034806, $null, -0.010657853425638181, 0.031468739706286171, [double]"3.7665825361947448e+18", [double]"1.275668646441314e+18", 0.031468739706286171, 0.23430532913982935),
    @(0, 2.7122120396162424, $null, -0.97942001697994874, 0.61367412211482841, 65535, 65535, 0.97942001697994874, 0.53369547192961131),
    @(0, 2.1421186466279774, $null, -0.099589484540936551, 0.22536832624391215, 65535, 65535, 0.22536832624391215, 0.07159329501494785)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = 18 + $i
    $vals = $data[$i]
    for ($col = 1; $col -le $vals.Length; $col++) {
        $v = $vals[$col - 1]
        $cell = $ws.Cells.Item($rowNum, $col)
        if ($null -ne $v) {
            $cell.Value = $v
        } else {
            # Materialize an explicit (but empty) cell, matching the
            # original sheet's blank "SettlingTime" column cells (e.g. C3:C5).
            $cell.Style = "Normal"
        }
    }
}
